$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The published template gained four new import columns ("form", "color",
# "solubility", "inventory label"). They were inserted as a block right
# after the existing "molecular mass (decoupled)"/"sum formula (decoupled)"
# pair (columns K:L), i.e. before what was column M ("flash point"),
# pushing every later column four places to the right.
$ws.Range("M1:P1").EntireColumn.Insert()

# Populate the header row for the newly inserted columns. (Column order on
# the sheet is color, solubility, inventory label, form.)
$ws.Range("M1").Value = "color"
$ws.Range("N1").Value = "solubility"
$ws.Range("O1").Value = "inventory label"
$ws.Range("P1").Value = "form"

# Match the width used for the other "highlighted" header columns in this
# block (column L / "sum formula (decoupled)" is ~22.86 characters wide).
$ws.Range("M1:P1").ColumnWidth = 22

# Restore the active selection to a cell within the newly added columns.
[void]$ws.Range("N5").Select()
